$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("I2").Value = 0.7964389134426562
$ws.Range("J2").Value = 0.7964389134426563
$ws.Range("M2").Value = 3.508776333333333
$ws.Range("N2").Value = 10.526329
$ws.Range("Q2").Value = 4.557558936103556
$ws.Range("R2").Value = 41.01803042493201
$ws.Range("S2").Value = 0.7964389134426562
$ws.Range("T2").Value = 0.7964389134426563

# Row 3 updates
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.3319853333333333
$ws.Range("H3").Value = 0.995956
$ws.Range("I3").Value = 0.2035610865573438
$ws.Range("J3").Value = 0.2035610865573438
$ws.Range("M3").Value = 3.508776333333333
$ws.Range("N3").Value = 10.526329
$ws.Range("Q3").Value = 1.164862280613778
$ws.Range("R3").Value = 10.483760525524
$ws.Range("S3").Value = 0.2035610865573438
$ws.Range("T3").Value = 0.2035610865573438
